$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 486.52173
$ws.Range("I107").Value = 468.57144
$ws.Range("J107").Value = 675
$ws.Range("K107").Value = 468.57144
$ws.Range("L107").Value = 675
$ws.Range("M107").Value = 1451.42856
$ws.Range("N107").Value = -4515

$ws.Range("H116").Value = 2368.5625
$ws.Range("I116").Value = 2129
$ws.Range("J116").Value = 2825.9092
$ws.Range("K116").Value = 2129
$ws.Range("L116").Value = 2825.9092
$ws.Range("M116").Value = 1313
$ws.Range("N116").Value = -9709.9092

$ws.Range("H125").Value = 2812.4546
$ws.Range("J125").Value = 2906.3333
$ws.Range("L125").Value = 26156.9997
$ws.Range("N125").Value = -31076.9997

$ws.Range("H132").Value = 9620774
$ws.Range("I132").Value = 12504587
$ws.Range("K132").Value = 37513761
$ws.Range("M132").Value = -37511231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 550
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = -1224

$ws.Range("H32").Value = 13335495
$ws.Range("I32").Value = 1909.9275
$ws.Range("K32").Value = 1909.9275
$ws.Range("M32").Value = -1622.9275

$ws.Range("H50").Value = 446.14285
$ws.Range("I50").Value = 348
$ws.Range("J50").Value = 462.5
$ws.Range("K50").Value = 348
$ws.Range("L50").Value = 462.5
$ws.Range("M50").Value = 366
$ws.Range("N50").Value = -1890.5

$ws.Range("H61").Value = 5557056.5
$ws.Range("I61").Value = 5748541
$ws.Range("K61").Value = 5748541
$ws.Range("M61").Value = -5748329

$ws.Range("H62").Value = 18800
$ws.Range("J62").Value = 18800
$ws.Range("L62").Value = 18800
$ws.Range("N62").Value = -20048

$ws.Range("H65").Value = 18800
$ws.Range("J65").Value = 18800
$ws.Range("L65").Value = 56400
$ws.Range("N65").Value = -62640

$ws.Range("H108").Value = 22925.334
$ws.Range("J108").Value = 22925.334
$ws.Range("L108").Value = 22925.334
$ws.Range("N108").Value = -30605.334

$ws.Range("H136").Value = 5557056.5
$ws.Range("I136").Value = 5748541
$ws.Range("K136").Value = 17245623
$ws.Range("M136").Value = -17243073

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 550
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = -1230

$ws.Range("H122").Value = 38585
$ws.Range("J122").Value = 38585
$ws.Range("L122").Value = 38585
$ws.Range("N122").Value = -48385

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 7410.7144
$ws.Range("I7").Value = 214
$ws.Range("J7").Value = 14607.429
$ws.Range("K7").Value = 214
$ws.Range("L7").Value = 14607.429
$ws.Range("M7").Value = -101
$ws.Range("N7").Value = -14833.429

$ws.Range("H35").Value = 3158.125
$ws.Range("I35").Value = 3158.125
$ws.Range("K35").Value = 3158.125
$ws.Range("M35").Value = -2864.125

$ws.Range("H132").Value = 11906386
$ws.Range("I132").Value = 1092.7142
$ws.Range("J132").Value = 47622268
$ws.Range("K132").Value = 3278.1426
$ws.Range("L132").Value = 142866804
$ws.Range("M132").Value = -748.1425999999997
$ws.Range("N132").Value = -142871864

$ws.Range("H134").Value = 14286566
$ws.Range("I134").Value = 960.5417
$ws.Range("J134").Value = 45455160
$ws.Range("K134").Value = 2881.6251
$ws.Range("L134").Value = 136365480
$ws.Range("M134").Value = -346.6251000000002
$ws.Range("N134").Value = -136370550

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 111306.89
$ws.Range("I4").Value = 111306.89
$ws.Range("K4").Value = 333920.67
$ws.Range("M4").Value = -333808.67

$ws.Range("H6").Value = 139.23077
$ws.Range("I6").Value = 91
$ws.Range("K6").Value = 273
$ws.Range("M6").Value = -160

$ws.Range("H16").Value = 1250.1666
$ws.Range("I16").Value = 667
$ws.Range("J16").Value = 1833.3334
$ws.Range("K16").Value = 2001
$ws.Range("L16").Value = 5500.0002
$ws.Range("M16").Value = -1828
$ws.Range("N16").Value = -5846.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H80").Value = 12501912
$ws.Range("I80").Value = 2060
$ws.Range("J80").Value = 33335000
$ws.Range("K80").Value = 2060
$ws.Range("L80").Value = 33335000
$ws.Range("M80").Value = -1062
$ws.Range("N80").Value = -33336996

$ws.Range("H83").Value = 12501912
$ws.Range("I83").Value = 2060
$ws.Range("J83").Value = 33335000
$ws.Range("K83").Value = 10300
$ws.Range("L83").Value = 166675000
$ws.Range("M83").Value = -5308
$ws.Range("N83").Value = -166684984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 153062350
$ws.Range("I16").Value = 14286286
$ws.Range("K16").Value = 14286286
$ws.Range("M16").Value = -14286116

$ws.Range("H132").Value = 6383.383
$ws.Range("I132").Value = 1555.963
$ws.Range("J132").Value = 12900.4
$ws.Range("K132").Value = 4667.889
$ws.Range("L132").Value = 38701.2
$ws.Range("M132").Value = -2137.889
$ws.Range("N132").Value = -43761.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 41823.5
$ws.Range("J12").Value = 55338
$ws.Range("L12").Value = 55338
$ws.Range("N12").Value = -55622

$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 20257.775
$ws.Range("I132").Value = 26765.46
$ws.Range("J132").Value = 6899.8945
$ws.Range("K132").Value = 80296.38
$ws.Range("L132").Value = 20699.6835
$ws.Range("M132").Value = -77766.38
$ws.Range("N132").Value = -25759.6835

$ws.Range("H136").Value = 31251954
$ws.Range("I136").Value = 43479230
$ws.Range("J136").Value = 4461.1113
$ws.Range("K136").Value = 130437690
$ws.Range("L136").Value = 13383.3339
$ws.Range("M136").Value = -130435140
$ws.Range("N136").Value = -18483.3339
